$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "login"

# Column widths (target XML widths 13 and 12.140625 characters; the engine
# snaps ColumnWidth to the nearest 1/6 character, so feed values that land
# in the correct rounding bucket)
$ws.Columns.Item(1).ColumnWidth = 12.17
$ws.Columns.Item(2).ColumnWidth = 11.33

# Column A
$ws.Range("A1").Value = "username"
$ws.Range("A2").Value = "selenium"

# Column B
$ws.Range("B1").Value = "password"
$ws.Range("B2").Value = "Selenium@123"

# Hyperlink on B2 (creates the Hyperlink cell style); no explicit display
# text is passed so the cell's own text is used as-is (no display= attr)
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Selenium@123")

# Highlight header row with yellow fill
$ws.Range("A1:B1").Interior.Color = 65535

$ws.Range("A1:B1").Select()
